# Implementing material ui design
# - Delete the empty "Video Games" sheet
# - Insert a new "Games" sheet right before "TV" (i.e. right after "Indie")
# - Populate "Games" with the new Persona 5 rows + the partially-filled
#   game/song backlog list
# - Activate the new "Games" sheet
# - Clear the stray fill style from Indie!H42:H43

$wb = $excel.ActiveWorkbook

# --- 1. Remove the old empty "Video Games" sheet -----------------------
$videoGames = $wb.Worksheets.Item("Video Games")
$videoGames.Delete()

# --- 2. Insert the new "Games" sheet just before "TV" -------------------
$tv = $wb.Worksheets.Item("TV")
$games = $wb.Worksheets.Add($tv)
$games.Name = "Games"

# --- 3. Header row --------------------------------------------------
$games.Range("A1").Value = "Name"
$games.Range("B1").Value = "Property"
$games.Range("C1").Value = "Difficulty"
$games.Range("D1").Value = "Song Name"
$games.Range("E1").Value = "Artist"
$games.Range("F1").Value = "Category"
$games.Range("G1").Value = "Location"
$games.Range("H1").Value = "Video Link"
$games.Range("I1").Value = "SQL"

# --- 4. Row 2 & 3 : Persona 5 Easy/Hard, fully filled in ---------------
$games.Range("A2").Value = "Persona 5 - Easy"
$games.Range("B2").Value = "Persona 5"
$games.Range("C2").Value = "Easy"
$games.Range("D2").Value = "Last Surprise"
$games.Range("E2").Value = "Shoji Meguro"
$games.Range("F2").Value = "Video Games"
$games.Range("G2").Formula = '="music/"&LOWER(F2)&"/"&A2'
$games.Range("H2").Value = "https://www.youtube.com/embed/Ec4YbVP9R-A?si=0e-P9iwxDCMkviKH"
$games.Range("I2").Formula = '="INSERT INTO songs (name, property, difficulty, song_name, arist, category, location, video_link) VALUES (''"&A2&"''"&", "&"''"&B2&"''"&", "&"''"&C2&"''"&", "&"''"&D2&"''"&", "&"''"&E2&"''"&", "&"''"&F2&"''"&", "&"''"&G2&"'', "&"''"&H2&"'');"'
$games.Range("I2").WrapText = $true
$games.Range("I2").VerticalAlignment = -4160

$games.Range("A3").Value = "Persona 5 - Hard"
$games.Range("B3").Value = "Persona 5"
$games.Range("C3").Value = "Hard"
$games.Range("D3").Value = "Last Surprise"
$games.Range("E3").Value = "Shoji Meguro"
$games.Range("F3").Value = "Video Games"
$games.Range("G3").Formula = '="music/"&LOWER(F3)&"/"&A3'
$games.Range("H3").Value = "https://www.youtube.com/embed/Ec4YbVP9R-A?si=0e-P9iwxDCMkviKH"
$games.Range("I3").Formula = '="INSERT INTO songs (name, property, difficulty, song_name, arist, category, location, video_link) VALUES (''"&A3&"''"&", "&"''"&B3&"''"&", "&"''"&C3&"''"&", "&"''"&D3&"''"&", "&"''"&E3&"''"&", "&"''"&F3&"''"&", "&"''"&G3&"'', "&"''"&H3&"'');"'
$games.Range("I3").WrapText = $true
$games.Range("I3").VerticalAlignment = -4160

# Row height for the two data-filled rows
$games.Range("A2:I3").RowHeight = 45

# --- 5. Rows 4-22: backlog list (only Property + Song Name known so far)
$backlog = @(
    @{ Row = 4;  B = "Metal Gear Rising: Revengeance"; D = "The Only Thing I Know For Real" },
    @{ Row = 5;  B = "Fire Emblem: Three Houses";       D = "Edge Of Dawn" },
    @{ Row = 6;  B = "Shin Megami Tensei III: Nocturne"; D = "Forced Battle" },
    @{ Row = 7;  B = "Persona 4" },
    @{ Row = 8;  B = "Persona 3";                       D = "Mass Destruction" },
    @{ Row = 9;  B = "Final Fantasy IX";                D = "The Place I'll Return To Someday" },
    @{ Row = 10; B = "Code Vein";                       D = "Requiem" },
    @{ Row = 11; B = "Dark Souls";                      D = "Ornstein & Smough" },
    @{ Row = 12; B = "Outer Wilds";                     D = "14.3 Billion Years"; H = "https://www.youtube.com/embed/XOrygf_iLhw?si=aLNfbxiTuwz1X_X0" },
    @{ Row = 13; B = "Shin Megami Tensei V" },
    @{ Row = 14; B = "Stardew Valley";                  D = "Spring (The Valley Comes Alive)" },
    @{ Row = 15; B = "Subnautica";                      D = "Abandon Ship" },
    @{ Row = 16; B = "Silent Hill 2";                   D = "Theme Of Laura" },
    @{ Row = 17; B = "Shadow Of The Colossus";          D = "Prologue" },
    @{ Row = 18; B = "Civilisation 4";                  D = "Baba Yetu" },
    @{ Row = 19; B = "The Elder Scrolls V: Skyrim";     D = "Dragonborn (Theme)" },
    @{ Row = 20; B = "Kingdom Hearts";                  D = "Dearly Beloved" },
    @{ Row = 21; B = "Final Fantasy X";                 D = "To Zanarkand" },
    @{ Row = 22; B = "Final Fantasy VII";                D = "One Winged Angel" }
)

foreach ($entry in $backlog) {
    $r = $entry.Row
    $games.Range("B$r").Value = $entry.B
    if ($entry.ContainsKey("D")) {
        $games.Range("D$r").Value = $entry.D
    }
    if ($entry.ContainsKey("H")) {
        $games.Range("H$r").Value = $entry.H
    }
}

# --- 6. Column widths / autofit (matches bestFit columns in the sheet) --
$games.Columns.Item("A:I").AutoFit()

# --- 7. Selection / view state, then make "Games" the active tab -------
$games.Range("D23").Select()
$excel.ActiveWindow.Zoom = 190
$games.Activate()

# --- 8. Clean up the leftover fill-only style on Indie!H42:H43 ---------
$indie = $wb.Worksheets.Item("Indie")
$indie.Range("H42:H43").Style = "Normal"
